$d = $word.ActiveDocument

# Locate the paragraph that holds the "<Proposal Description>" placeholder
# text and insert a new, blank paragraph (same "Keybody" / italic style)
# immediately before it, giving the consultation proposal an extra blank
# line above it.
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*<Proposal Description>*") {
        $r = $para.Range.Duplicate
        $r.Collapse(1)
        $r.InsertParagraphBefore()
        break
    }
}

$d.Save()
